# TEST GROUP DATA CLEANING AND EDA
# - Convert the control_group data range into a native Excel Table
#   (ListObject), matching the column headers already in row 1.
# - Resize the data columns (no longer "best fit", now explicit custom
#   widths to comfortably show the table's filter-drop-down buttons).
# - Move the active selection off the whole-sheet default onto O22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Turn A1:J31 into an Excel Table named "control_group" -----------------
$dataRange = $ws.Range("A1:J31")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "control_group"
$tbl.TableStyle = "TableStyleMedium2"

# --- Widen the columns now that they host table header/filter buttons ------
$ws.Columns.Item(1).ColumnWidth = 16.73697916666667
$ws.Columns.Item(3).ColumnWidth = 13.451822916666666
$ws.Columns.Item(4).ColumnWidth = 16.73697916666667
$ws.Columns.Item(5).ColumnWidth = 7.592447916666667
$ws.Columns.Item(6).ColumnWidth = 19.022135416666668
$ws.Columns.Item(7).ColumnWidth = 13.877604166666666
$ws.Columns.Item(8).ColumnWidth = 18.166666666666668
$ws.Columns.Item(9).ColumnWidth = 15.877604166666666
$ws.Columns.Item(10).ColumnWidth = 14.022135416666666

# --- Move the selection to O22 (instead of the whole-sheet default) --------
$ws.Range("O22").Select() | Out-Null
